$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (Date serial in column D, Volumen in column J,
# and the shared price in columns K, L, M, P).
$rows = @{
    2  = @{ D = 44365; J = 55; P = 5000 }
    3  = @{ D = 44312; J = 50; P = 4000 }
    4  = @{ D = 44259; J = 30; P = 4000 }
    5  = @{ D = 44280; J = 55; P = 4000 }
    6  = @{ D = 44390; J = 55; P = 6000 }
    7  = @{ D = 44509; J = 20; P = 4000 }
    8  = @{ D = 44316; J = 20; P = 4000 }
    9  = @{ D = 44313; J = 20; P = 4000 }
    12 = @{ D = 44315; J = 40; P = 4000 }
    13 = @{ D = 44504; J = 55; P = 4000 }
    14 = @{ D = 44508; J = 30; P = 4000 }
    15 = @{ D = 44291; J = 35; P = 4000 }
    16 = @{ D = 44301; J = 40; P = 3000 }
    18 = @{ D = 44497; J = 20; P = 4000 }
    19 = @{ D = 44176; J = 10; P = 4000 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value = $vals.D
    $ws.Cells.Item($r, 10).Value = $vals.J
    $ws.Cells.Item($r, 11).Value = $vals.P
    $ws.Cells.Item($r, 12).Value = $vals.P
    $ws.Cells.Item($r, 13).Value = $vals.P
    $ws.Cells.Item($r, 16).Value = $vals.P
}
